# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price (D) and Volume(1h) (E) columns are stored as plain text in the source
# sheet, so numeric-looking prices get `NumberFormat = '@'` first to stop
# Excel's COM layer from auto-coercing them into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.848.91'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '2.687.87'
$ws.Range('E3').Value = '  +4.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '513.80'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.03'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  +1.59%  '
$ws.Range('D9').Value = '2.684.55'
$ws.Range('E9').Value = '  +4.29%  '
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('E11').Value = '  +4.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.335'
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').Value = '3.121.99'
$ws.Range('E14').Value = '  +3.51%  '
$ws.Range('D15').Value = '58.829.01'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.95'
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '2.673.60'
$ws.Range('E18').Value = '  +4.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '347.74'
$ws.Range('E19').Value = '  +4.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.53'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  +3.14%  '
$ws.Range('E22').Value = '  +2.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.86'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('E25').Value = '  +2.75%  '
$ws.Range('D26').Value = '2.773.31'
$ws.Range('E26').Value = '  +3.86%  '
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('E28').Value = '  +1.26%  '
$ws.Range('E29').Value = '  +3.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.23'
$ws.Range('E30').Value = '  +4.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.996'
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.38'
$ws.Range('E32').Value = '  +8.88%  '
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('E34').Value = '  +1.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '149.38'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('E36').Value = '  +12.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.01'
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('E38').Value = '  +3.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.73'
$ws.Range('E39').Value = '  +2.16%  '
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.68'
$ws.Range('E41').Value = '  +4.15%  '
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '278.26'
$ws.Range('E44').Value = '  -3.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0979'
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.59'
$ws.Range('E47').Value = '  +4.07%  '
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.72'
$ws.Range('E49').Value = '  +4.33%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0229'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.27'
$ws.Range('E51').Value = '  -0.50%  '
